$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Group"

$ws.Range("E2").Value = 309.8567
$ws.Range("E3").Value = 368.4094
$ws.Range("E4").Value = 0.7143
$ws.Range("E6").Value = 0.178
